# Update cryptocurrency price/volume data on the active sheet.
# Values that parse as plain numbers use a leading apostrophe (quote-prefix)
# via .Formula so Excel stores them as text (matching the source data),
# instead of silently converting them to numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.824.00"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").Value = "2.450.59"
$ws.Range("E3").Value = "  +1.12%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Formula = "'570.45"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Formula = "'146.64"
$ws.Range("E6").Value = "  +2.05%  "
$ws.Range("D7").Formula = "'1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Formula = "'0.534"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "2.452.91"
$ws.Range("E9").Value = "  +1.25%  "
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Formula = "'5.24"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Formula = "'0.356"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").Formula = "'26.87"
$ws.Range("E14").Value = "  +2.68%  "
$ws.Range("E15").Value = "  +3.29%  "
$ws.Range("D16").Value = "2.907.92"
$ws.Range("D17").Value = "63.062.03"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "2.457.14"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("D19").Formula = "'11.33"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").Formula = "'7.24"
$ws.Range("E20").Value = "  +6.24%  "
$ws.Range("D21").Formula = "'323.17"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  +0.67%  "
$ws.Range("D23").Formula = "'1.95"
$ws.Range("E23").Value = "  +13.13%  "
$ws.Range("D24").Formula = "'1.00"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Formula = "'66.17"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").Formula = "'621.23"
$ws.Range("E26").Value = "  +11.31%  "
$ws.Range("D27").Formula = "'8.60"
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  +9.42%  "
$ws.Range("D29").Value = "2.567.81"
$ws.Range("D30").Formula = "'0.995"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("E31").Value = "  +6.19%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("D33").Formula = "'0.141"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("E34").Value = "  +1.86%  "
$ws.Range("E35").Value = "  +6.68%  "
$ws.Range("D36").Formula = "'1.48"
$ws.Range("E36").Value = "  -1.53%  "
$ws.Range("D37").Formula = "'0.999"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Formula = "'0.382"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").Formula = "'5.40"
$ws.Range("E39").Value = "  -1.38%  "
$ws.Range("D40").Formula = "'18.67"
$ws.Range("E40").Value = "  -0.04%  "
$ws.Range("D41").Formula = "'145.21"
$ws.Range("E41").Value = "  -4.43%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("E43").Value = "  +16.20%  "
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("D45").Formula = "'147.29"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Formula = "'3.72"
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Formula = "'20.70"
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Formula = "'0.0537"
$ws.Range("E48").Value = "  +1.56%  "
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").Formula = "'0.0234"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").Formula = "'0.0918"
$ws.Range("E51").Value = "  -0.55%  "
